$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '58.577.72'
$ws.Range("E2").Value = '  -0.66%  '
$ws.Range("D3").Value = '2.487.91'
$ws.Range("E3").Value = '  -0.48%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '527.75'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.88%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '134.15'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -3.25%  '
$ws.Range("E7").Value = '  +0.34%  '
$ws.Range("E8").Value = '  -0.77%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.100'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.96%  '
$ws.Range("E10").Value = '  -2.14%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.39'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.41%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.344'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.77%  '
$ws.Range("D13").Value = '2.929.40'
$ws.Range("E13").Value = '  -0.74%  '
$ws.Range("D14").Value = '58.472.88'
$ws.Range("E14").Value = '  -0.66%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '22.50'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -3.19%  '
$ws.Range("E16").Value = '  -1.85%  '
$ws.Range("D17").Value = '2.484.90'
$ws.Range("E17").Value = '  -1.66%  '
$ws.Range("E18").Value = '  -1.34%  '
$ws.Range("E19").Value = '  -1.73%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '321.05'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.61%  '
$ws.Range("E21").Value = '  +0.10%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.83'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.16%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '64.46'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.63%  '
$ws.Range("E24").Value = '  -2.48%  '
$ws.Range("B25").Value = 'Binance-PegBSC-USD'
$ws.Range("C25").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.00'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.06%  '
$ws.Range("B26").Value = 'Kaspa'
$ws.Range("C26").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.162'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.34%  '
$ws.Range("E27").Value = '  -2.63%  '
$ws.Range("D28").Value = '0.0₃0754'
$ws.Range("E28").Value = '  -3.22%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '6.44'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -4.95%  '
$ws.Range("E30").Value = '  -3.26%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '165.93'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -2.11%  '
$ws.Range("E32").Value = '  -5.11%  '
$ws.Range("E33").Value = '  -0.04%  '
$ws.Range("E34").Value = '  +0.39%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '18.29'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.42%  '
$ws.Range("E36").Value = '  -8.83%  '
$ws.Range("E37").Value = '  -3.18%  '
$ws.Range("E38").Value = '  -4.05%  '
$ws.Range("E39").Value = '  -3.17%  '
$ws.Range("E40").Value = '  -3.12%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '276.79'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -2.58%  '
$ws.Range("E42").Value = '  -6.60%  '
$ws.Range("E43").Value = '  -1.30%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '127.49'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.54%  '
$ws.Range("E45").Value = '  -1.82%  '
$ws.Range("E46").Value = '  -2.97%  '
$ws.Range("E47").Value = '  -2.63%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '17.27'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.70%  '
$ws.Range("D49").Value = '1.742.75'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.979'
$ws.Range("D50").Style = "Normal"
$ws.Range("E51").Value = '  -1.77%  '
